# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels AD1:AF1 ---
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Match the header style used by the rest of row 1 (bold, centered, bordered)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-47): same season record (83 wins, 79 losses, 0 ties) for every player ---
For ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 83
    $ws.Cells.Item($r, 31).Value2 = 79
    $ws.Cells.Item($r, 32).Value2 = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-47"
